# Applies the F/G column value updates described by the commit diff
# (gh-pages data refresh at commit 456a3b4) across all four worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 41190
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 16
$ws.Range("F5").Value = 9189
$ws.Range("F6").Value = 191
$ws.Range("F7").Value = 717
$ws.Range("F8").Value = 828
$ws.Range("F9").Value = 685
$ws.Range("F10").Value = 184
$ws.Range("F12").Value = 270
$ws.Range("F13").Value = 825
$ws.Range("F14").Value = 74
$ws.Range("F15").Value = 110
$ws.Range("F17").Value = 285
$ws.Range("F18").Value = 1299
$ws.Range("F20").Value = 555
$ws.Range("F21").Value = 663
$ws.Range("F22").Value = 440
$ws.Range("F23").Value = 649
$ws.Range("F24").Value = 692
$ws.Range("F27").Value = 53
$ws.Range("F28").Value = 452
$ws.Range("F29").Value = 484
$ws.Range("F32").Value = 899
$ws.Range("F33").Value = 418
$ws.Range("F34").Value = 73
$ws.Range("F35").Value = 202
$ws.Range("F36").Value = 128
$ws.Range("F37").Value = 332
$ws.Range("F38").Value = 1174
$ws.Range("F39").Value = 268
$ws.Range("F41").Value = 1181
$ws.Range("F44").Value = 10
$ws.Range("F46").Value = 33

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 33
$ws.Range("F5").Value = 4421
$ws.Range("F11").Value = 114

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1967
$ws.Range("F3").Value = 471
$ws.Range("F4").Value = 321

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1967
$ws.Range("F3").Value = 471
$ws.Range("F4").Value = 41190
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 33
$ws.Range("F11").Value = 9189
$ws.Range("F12").Value = 191
$ws.Range("F13").Value = 718
$ws.Range("F14").Value = 718
$ws.Range("F16").Value = 321
$ws.Range("F17").Value = 828
$ws.Range("F18").Value = 114
$ws.Range("F19").Value = 270
$ws.Range("F20").Value = 825
$ws.Range("F21").Value = 74
$ws.Range("F22").Value = 110
$ws.Range("F23").Value = 285
$ws.Range("F24").Value = 1299
$ws.Range("F25").Value = 555
$ws.Range("F26").Value = 440
$ws.Range("F27").Value = 649
$ws.Range("F28").Value = 692
$ws.Range("F30").Value = 53
$ws.Range("F31").Value = 452
$ws.Range("F34").Value = 484
$ws.Range("F37").Value = 899
$ws.Range("F39").Value = 418
$ws.Range("F40").Value = 73
$ws.Range("F41").Value = 128
$ws.Range("F42").Value = 332
$ws.Range("F43").Value = 268
$ws.Range("F44").Value = 1181
$ws.Range("F48").Value = 10
